$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 470.9091
$ws.Range("I6").Value = 175.66667
$ws.Range("J6").Value = 1799.5
$ws.Range("K6").Value = 527.00001
$ws.Range("L6").Value = 5398.5
$ws.Range("M6").Value = -415.00001
$ws.Range("N6").Value = -5622.5

$ws.Range("H9").Value = 107.083336
$ws.Range("I9").Value = 98.40000000000001
$ws.Range("J9").Value = 150.5
$ws.Range("K9").Value = 98.40000000000001
$ws.Range("L9").Value = 150.5
$ws.Range("M9").Value = 70.59999999999999
$ws.Range("N9").Value = -488.5

$ws.Range("H17").Value = 469.48718
$ws.Range("J17").Value = 348.60526
$ws.Range("L17").Value = 1045.81578
$ws.Range("N17").Value = -1381.81578

$ws.Range("H92").Value = 1308.2
$ws.Range("I92").Value = 1308.2
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 1308.2
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -60.20000000000005
$ws.Range("N92").ClearContents()

$ws.Range("H137").Value = 2528.111
$ws.Range("J137").Value = 2704.32
$ws.Range("L137").Value = 8112.960000000001
$ws.Range("N137").Value = -13212.96

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12249.36
$ws.Range("I32").Value = 13419.075
$ws.Range("J32").Value = 4500
$ws.Range("K32").Value = 13419.075
$ws.Range("L32").Value = 4500
$ws.Range("M32").Value = -13132.075
$ws.Range("N32").Value = -5074

$ws.Range("H74").Value = 1166.8511
$ws.Range("I74").Value = 1031.0555
$ws.Range("J74").Value = 1611.2727
$ws.Range("K74").Value = 1031.0555
$ws.Range("L74").Value = 1611.2727
$ws.Range("M74").Value = -157.0554999999999
$ws.Range("N74").Value = -3359.2727

$ws.Range("H77").Value = 1166.8511
$ws.Range("I77").Value = 1031.0555
$ws.Range("J77").Value = 1611.2727
$ws.Range("K77").Value = 5155.2775
$ws.Range("L77").Value = 8056.363499999999
$ws.Range("M77").Value = -787.2775000000001
$ws.Range("N77").Value = -16792.3635

$ws.Range("H110").Value = 3099.4546
$ws.Range("I110").Value = 2399.1428
$ws.Range("J110").Value = 4325
$ws.Range("K110").Value = 2399.1428
$ws.Range("L110").Value = 4325
$ws.Range("M110").Value = -354.1428000000001
$ws.Range("N110").Value = -8415

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2438.6428
$ws.Range("I99").Value = 2241.6667
$ws.Range("J99").Value = 2586.375
$ws.Range("K99").Value = 2241.6667
$ws.Range("L99").Value = 2586.375
$ws.Range("M99").Value = -743.6667000000002
$ws.Range("N99").Value = -5582.375

$ws.Range("H105").Value = 3599.4
$ws.Range("I105").Value = 3499.5
$ws.Range("K105").Value = 3499.5
$ws.Range("M105").Value = -1752.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 7642.3335
$ws.Range("I4").Value = 99.8
$ws.Range("J4").Value = 9999.375
$ws.Range("K4").Value = 99.8
$ws.Range("L4").Value = 9999.375
$ws.Range("M4").Value = 12.2
$ws.Range("N4").Value = -10223.375

$ws.Range("H31").Value = 1981.8448
$ws.Range("I31").Value = 2648.4348
$ws.Range("J31").Value = 1543.8
$ws.Range("K31").Value = 2648.4348
$ws.Range("L31").Value = 1543.8
$ws.Range("M31").Value = -2353.4348
$ws.Range("N31").Value = -2133.8

$ws.Range("H34").Value = 1981.8448
$ws.Range("I34").Value = 2648.4348
$ws.Range("J34").Value = 1543.8
$ws.Range("K34").Value = 2648.4348
$ws.Range("L34").Value = 1543.8
$ws.Range("M34").Value = -2446.4348
$ws.Range("N34").Value = -1947.8

$ws.Range("H35").Value = 651.125
$ws.Range("I35").Value = 651.125
$ws.Range("K35").Value = 651.125
$ws.Range("M35").Value = -357.125

$ws.Range("H96").Value = 28467.75
$ws.Range("J96").Value = 28467.75
$ws.Range("L96").Value = 28467.75
$ws.Range("N96").Value = -33959.75

$ws.Range("H134").Value = 2650.8948
$ws.Range("I134").Value = 1688.5
$ws.Range("J134").Value = 5345.6
$ws.Range("K134").Value = 5065.5
$ws.Range("L134").Value = 16036.8
$ws.Range("M134").Value = -2530.5
$ws.Range("N134").Value = -21106.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 1355
$ws.Range("J29").Value = 1226
$ws.Range("L29").Value = 3678
$ws.Range("N29").Value = -4232

$ws.Range("H68").Value = 807.0102000000001
$ws.Range("I68").Value = 631.9355
$ws.Range("J68").Value = 1108.5278
$ws.Range("K68").Value = 1895.8065
$ws.Range("L68").Value = 3325.5834
$ws.Range("M68").Value = -1084.8065
$ws.Range("N68").Value = -4947.5834

$ws.Range("H71").Value = 807.0102000000001
$ws.Range("I71").Value = 631.9355
$ws.Range("J71").Value = 1108.5278
$ws.Range("K71").Value = 5687.4195
$ws.Range("L71").Value = 9976.7502
$ws.Range("M71").Value = -1631.4195
$ws.Range("N71").Value = -18088.7502

$ws.Range("H107").Value = 1111.141
$ws.Range("I107").Value = 1201.3529
$ws.Range("J107").Value = 1041.4318
$ws.Range("K107").Value = 3604.0587
$ws.Range("L107").Value = 3124.2954
$ws.Range("M107").Value = -1684.0587
$ws.Range("N107").Value = -6964.2954

$ws.Range("H131").Value = 13891408
$ws.Range("I131").Value = 589.5454999999999
$ws.Range("K131").Value = 1768.6365
$ws.Range("M131").Value = 3271.3635

$ws.Range("H133").Value = 4462.25
$ws.Range("I133").Value = 1693.625
$ws.Range("J133").Value = 9999.5
$ws.Range("K133").Value = 5080.875
$ws.Range("L133").Value = 29998.5
$ws.Range("M133").Value = -20.875
$ws.Range("N133").Value = -40118.5

$ws.Range("H136").Value = 2125.7144
$ws.Range("I136").Value = 1376
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 4128
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = 972
$ws.Range("N136").Value = -22200

$ws.Range("H139").Value = 1938.52
$ws.Range("I139").Value = 1603.3158
$ws.Range("J139").Value = 3000
$ws.Range("K139").Value = 4809.9474
$ws.Range("L139").Value = 9000
$ws.Range("M139").Value = 330.0526
$ws.Range("N139").Value = -19280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 364850.2
$ws.Range("I70").Value = 534663.75
$ws.Range("J70").Value = 6354.8887
$ws.Range("K70").Value = 534663.75
$ws.Range("L70").Value = 6354.8887
$ws.Range("M70").Value = -534393.75
$ws.Range("N70").Value = -6894.8887

$ws.Range("H73").Value = 364850.2
$ws.Range("I73").Value = 534663.75
$ws.Range("J73").Value = 6354.8887
$ws.Range("K73").Value = 534663.75
$ws.Range("L73").Value = 6354.8887
$ws.Range("M73").Value = -533727.75
$ws.Range("N73").Value = -8226.8887

$ws.Range("H92").Value = 40251
$ws.Range("J92").Value = 40251
$ws.Range("L92").Value = 40251
$ws.Range("N92").Value = -43995

$ws.Range("H113").Value = 2444.0833
$ws.Range("I113").Value = 1500
$ws.Range("J113").Value = 2632.9
$ws.Range("K113").Value = 1500
$ws.Range("L113").Value = 2632.9
$ws.Range("M113").Value = 670
$ws.Range("N113").Value = -6972.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 11830.737
$ws.Range("I61").Value = 13686.5
$ws.Range("J61").Value = 1933.3334
$ws.Range("K61").Value = 13686.5
$ws.Range("L61").Value = 1933.3334
$ws.Range("M61").Value = -13484.5
$ws.Range("N61").Value = -2337.3334

$ws.Range("H93").Value = 1445.8182
$ws.Range("I93").Value = 1228.5714
$ws.Range("K93").Value = 1228.5714
$ws.Range("M93").Value = 19.42859999999996

$ws.Range("H94").Value = 24440
$ws.Range("J94").Value = 24440
$ws.Range("L94").Value = 24440
$ws.Range("N94").Value = -25792

$ws.Range("H113").Value = 11830.737
$ws.Range("I113").Value = 13686.5
$ws.Range("J113").Value = 1933.3334
$ws.Range("K113").Value = 13686.5
$ws.Range("L113").Value = 1933.3334
$ws.Range("M113").Value = -11516.5
$ws.Range("N113").Value = -6273.3334

$ws.Range("H132").Value = 2762.55
$ws.Range("I132").Value = 2627.2808
$ws.Range("J132").Value = 5332.6665
$ws.Range("K132").Value = 7881.8424
$ws.Range("L132").Value = 15997.9995
$ws.Range("M132").Value = -5351.8424
$ws.Range("N132").Value = -21057.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 23999
$ws.Range("J16").Value = 23999
$ws.Range("L16").Value = 23999
$ws.Range("N16").Value = -24583

$ws.Range("H19").Value = 1000
$ws.Range("J19").Value = 1000
$ws.Range("L19").Value = 1000
$ws.Range("N19").Value = -1348

$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

$ws.Range("H126").Value = 7402.0454
$ws.Range("I126").Value = 9677.75
$ws.Range("J126").Value = 1333.5
$ws.Range("K126").Value = 29033.25
$ws.Range("L126").Value = 4000.5
$ws.Range("M126").Value = -26563.25
$ws.Range("N126").Value = -8940.5

$ws.Range("H132").Value = 1657.3273
$ws.Range("I132").Value = 1508.921
$ws.Range("J132").Value = 1989.0588
$ws.Range("K132").Value = 4526.763
$ws.Range("L132").Value = 5967.1764
$ws.Range("M132").Value = -1996.763
$ws.Range("N132").Value = -11027.1764

